$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.961.11"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "2.994.09"
$ws.Range("E3").Value = "  +2.38%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'353.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").Value = "'106.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.71%  "
$ws.Range("E7").Value = "  -2.42%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  -3.70%  "
$ws.Range("D10").Value = "'37.87"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.68%  "
$ws.Range("E11").Value = "  +2.60%  "
$ws.Range("D12").Value = "'0.0853"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.47%  "
$ws.Range("D13").Value = "'18.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.85%  "
$ws.Range("D14").Value = "3.467.50"
$ws.Range("E14").Value = "  +2.43%  "
$ws.Range("D15").Value = "'7.56"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.12%  "
$ws.Range("D16").Value = "2.982.52"
$ws.Range("E16").Value = "  +2.26%  "
$ws.Range("D17").Value = "'1.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.37%  "
$ws.Range("D18").Value = "51.912.01"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").Value = "'3.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.13%  "
$ws.Range("D20").Value = "'7.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.35%  "
$ws.Range("E21").Value = "  -3.80%  "
$ws.Range("D22").Value = "0.0₃0967"
$ws.Range("E22").Value = "  -1.50%  "
$ws.Range("D23").Value = "'68.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.79%  "
$ws.Range("D24").Value = "'262.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.92%  "
$ws.Range("E25").Value = "  -3.85%  "
$ws.Range("E26").Value = "  -2.46%  "
$ws.Range("D27").Value = "'26.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.97%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("E29").Value = "  +0.90%  "
$ws.Range("D30").Value = "'0.108"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.56%  "
$ws.Range("D31").Value = "'6.36"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.92%  "
$ws.Range("D32").Value = "'10.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.06%  "
$ws.Range("D33").Value = "'35.90"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.17%  "
$ws.Range("D34").Value = "'2.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +13.04%  "
$ws.Range("D35").Value = "'51.20"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.93%  "
$ws.Range("D36").Value = "'0.0432"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.99%  "
$ws.Range("D37").Value = "'0.998"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("E38").Value = "  +1.70%  "
$ws.Range("D39").Value = "'2.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.58%  "
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").Value = "'17.44"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.25%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'1.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.95%  "
$ws.Range("E42").Value = "  -3.03%  "
$ws.Range("D43").Value = "'23.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.04%  "
$ws.Range("D44").Value = "'123.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.81%  "
$ws.Range("D45").Value = "'2.18"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("D46").Value = "2.122.34"
$ws.Range("E46").Value = "  -0.80%  "
$ws.Range("E47").Value = "  -4.30%  "
$ws.Range("E48").Value = "  -7.38%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "3.295.09"
$ws.Range("E49").Value = "  +2.60%  "
$ws.Range("B50").Value = "TheGraph"
$ws.Range("C50").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D50").Value = "'0.241"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.64%  "
$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D51").Value = "'0.0333"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.71%  "
